# Insert a new data row at row 117 (pushing the existing rows 117..190 down
# to 118..191) and populate the new row with the new Camote / "1a nueva(o)"
# reading dated 2022-02-11 (serial 44603).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("117:117").Insert()

$ws.Range("A117").Value = 5
$ws.Range("B117").Value = "Macroferia Regional de Talca"
$ws.Range("C117").Value = "Maule"
$ws.Range("D117").Value = 44603
$ws.Range("E117").Value = 7
$ws.Range("F117").Value = 100112045
$ws.Range("G117").Value = "Zapallo"
$ws.Range("H117").Value = "Camote"
$ws.Range("I117").Value = "1a nueva(o)"
$ws.Range("J117").Value = 900
$ws.Range("K117").Value = 300
$ws.Range("L117").Value = 300
$ws.Range("M117").Value = 300
$ws.Range("N117").Value = "$/kilo (volumen en unidades)"
$ws.Range("O117").Value = "Región del Maule"
$ws.Range("P117").Value = 300
$ws.Range("Q117").Value = 1
$ws.Range("R117").Value = "Hortaliza"
